$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "29.166.52"
$ws.Range("E2").Value = "  +3.02%  "

Set-TextValue $ws.Range("D3") "1.578.68"
$ws.Range("E3").Value = "  +1.76%  "

Set-TextValue $ws.Range("D4") "0.998"
$ws.Range("E4").Value = "  -0.19%  "

Set-TextValue $ws.Range("D5") "212.30"
$ws.Range("E5").Value = "  +1.10%  "

$ws.Range("E6").Value = "  +5.73%  "

Set-TextValue $ws.Range("D8") "26.07"
$ws.Range("E8").Value = "  +9.76%  "

Set-TextValue $ws.Range("D9") "0.248"
$ws.Range("E9").Value = "  +2.30%  "

Set-TextValue $ws.Range("D10") "0.0593"
$ws.Range("E10").Value = "  +1.68%  "

Set-TextValue $ws.Range("D11") "0.0906"
$ws.Range("E11").Value = "  +1.75%  "

Set-TextValue $ws.Range("D12") "1.806.73"
$ws.Range("E12").Value = "  +1.86%  "

Set-TextValue $ws.Range("D13") "1.550.82"
$ws.Range("E13").Value = "  -0.07%  "

Set-TextValue $ws.Range("D14") "29.194.08"
$ws.Range("E14").Value = "  +3.21%  "

Set-TextValue $ws.Range("D15") "0.522"
$ws.Range("E15").Value = "  +2.43%  "

Set-TextValue $ws.Range("D16") "3.70"
$ws.Range("E16").Value = "  +2.18%  "

Set-TextValue $ws.Range("D17") "62.18"
$ws.Range("E17").Value = "  +2.67%  "

Set-TextValue $ws.Range("D18") "236.67"
$ws.Range("E18").Value = "  +3.89%  "

Set-TextValue $ws.Range("D19") "7.43"
$ws.Range("E19").Value = "  +1.41%  "

Set-TextValue $ws.Range("D20") "0.0₃0689"
$ws.Range("E20").Value = "  +2.24%  "

Set-TextValue $ws.Range("D22") "3.98"
$ws.Range("E22").Value = "  +1.76%  "

Set-TextValue $ws.Range("D23") "9.16"
$ws.Range("E23").Value = "  +2.82%  "

Set-TextValue $ws.Range("D24") "2.07"
$ws.Range("E24").Value = "  +2.55%  "

Set-TextValue $ws.Range("D25") "153.72"
$ws.Range("E25").Value = "  +1.77%  "

Set-TextValue $ws.Range("D26") "15.14"
$ws.Range("E26").Value = "  +2.67%  "

Set-TextValue $ws.Range("D27") "0.107"
$ws.Range("E27").Value = "  +4.28%  "

Set-TextValue $ws.Range("D28") "6.34"
$ws.Range("E28").Value = "  +1.54%  "

$ws.Range("E29").Value = "  -0.12%  "

$ws.Range("E30").Value = "  +0.27%  "

Set-TextValue $ws.Range("D31") "1.06"
$ws.Range("E31").Value = "  +0.16%  "

$ws.Range("E32").Value = "  +1.55%  "

Set-TextValue $ws.Range("D33") "1.419.55"
$ws.Range("E33").Value = "  +2.35%  "

Set-TextValue $ws.Range("D34") "3.06"
$ws.Range("E34").Value = "  +1.53%  "

Set-TextValue $ws.Range("D35") "1.04"
$ws.Range("E35").Value = "  -2.71%  "

$ws.Range("E36").Value = "  +1.93%  "

Set-TextValue $ws.Range("D37") "2.75"
$ws.Range("E37").Value = "  +6.21%  "

$ws.Range("E38").Value = "  -1.76%  "

Set-TextValue $ws.Range("D40") "0.528"
$ws.Range("E40").Value = "  +3.18%  "

Set-TextValue $ws.Range("D41") "1.96"
$ws.Range("E41").Value = "  +2.44%  "

Set-TextValue $ws.Range("D42") "53.22"
$ws.Range("E42").Value = "  +24.22%  "

$ws.Range("E43").Value = "  -0.16%  "

Set-TextValue $ws.Range("D44") "0.789"
$ws.Range("E44").Value = "  +1.48%  "

Set-TextValue $ws.Range("D45") "0.0472"
$ws.Range("E45").Value = "  +1.67%  "

Set-TextValue $ws.Range("D46") "64.44"
$ws.Range("E46").Value = "  +4.32%  "

Set-TextValue $ws.Range("D47") "5.33"
$ws.Range("E47").Value = "  -0.31%  "

Set-TextValue $ws.Range("D48") "1.717.68"
$ws.Range("E48").Value = "  +1.86%  "

Set-TextValue $ws.Range("D49") "0.843"
$ws.Range("E49").Value = "  -6.93%  "

Set-TextValue $ws.Range("D50") "85.33"
$ws.Range("E50").Value = "  -0.33%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D51") "0.0516"
$ws.Range("E51").Value = "  +1.54%  "
